# Updated symbol list on Sun Dec 18 20:47:45 UTC 2022 with GitHub Actions
#
# Refresh the cryptocurrency price snapshot: most rows only get a new
# "Price" (column D) figure, while rows 16/17 swap their Coin/Link/Price/
# Volume data (One <-> CoinExToken changed rank) and row 43's Volume label
# loses its "Bestin24h" suffix alongside its new price.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$CellRef,
        [string]$Text
    )
    # Column D prices are stored as plain text even though they look like
    # numbers (e.g. "249.57"). A bare assignment would let Excel coerce the
    # string to a numeric cell, so numeric-looking values are entered with a
    # leading apostrophe - exactly like a user forcing text entry - to keep
    # them as text without touching the cell's number format. The leading
    # apostrophe itself only marks the quote-prefix on the cell's style, so
    # the style is reset back to Normal afterwards to leave formatting as it
    # was originally (plain General-formatted text cell).
    $needsTextForce = $Text -match '^[+-]?\d+(\.\d+)?$'
    $rng = $ws.Range($CellRef)
    if ($needsTextForce) {
        $rng.Value = "'" + $Text
        $rng.Style = "Normal"
    } else {
        $rng.Value = $Text
    }
}

# --- Simple price-only updates (column D) ---
Set-TextValue "D2"  "249.57"
Set-TextValue "D3"  "21.95"
Set-TextValue "D4"  "5.536"
Set-TextValue "D6"  "6.458"
Set-TextValue "D7"  "0.8011"
Set-TextValue "D8"  "1.038"
Set-TextValue "D10" "0.07302"
Set-TextValue "D11" "0.03096"
Set-TextValue "D12" "0.02912"
Set-TextValue "D13" "0.09271"
Set-TextValue "D14" "0.001670"
Set-TextValue "D15" "3.215"

# --- Rows 16 & 17 swap: "One" moves up to rank 15, "CoinExToken" to rank 16 ---
Set-TextValue "B16" "One"
Set-TextValue "C16" "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextValue "D16" "0.001296"
Set-TextValue "E16" "15OneONEBestin24h"

Set-TextValue "B17" "CoinExToken"
Set-TextValue "C17" "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
Set-TextValue "D17" "0.04740"
Set-TextValue "E17" "16CoinExTokenCET"

# --- More price-only updates (column D) ---
Set-TextValue "D18" "0.006450"
Set-TextValue "D20" "0.001052"
Set-TextValue "D22" "3.978"
Set-TextValue "D23" "3.374"
Set-TextValue "D24" "2.114"
Set-TextValue "D25" "0.3269"
Set-TextValue "D27" "0.0003300"
Set-TextValue "D40" "0.04148"
Set-TextValue "D41" "0.006901"

# --- Row 43: new price and the Volume label drops its "Bestin24h" suffix ---
Set-TextValue "D43" "0.002970"
Set-TextValue "E43" "42CEJICEJI"

# --- More price-only updates (column D) ---
Set-TextValue "D45" "0.00005637"
Set-TextValue "D48" "0.01636"
